$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# row 74
$ws.Range("H74").Value = 6582958
$ws.Range("I74").Value = 3640.4
$ws.Range("K74").Value = 3640.4
$ws.Range("M74").Value = -2704.4

# row 76
$ws.Range("H76").Value = 6175717.5
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 27779528
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 27779528
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -27780158

# row 77
$ws.Range("H77").Value = 6582958
$ws.Range("I77").Value = 3640.4
$ws.Range("K77").Value = 18202
$ws.Range("M77").Value = -13522

# row 79
$ws.Range("H79").Value = 6175717.5
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 27779528
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 27779528
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -27781712

# row 98
$ws.Range("H98").Value = 1059.6
$ws.Range("I98").Value = 1000
$ws.Range("J98").Value = 1298
$ws.Range("K98").Value = 1000
$ws.Range("L98").Value = 1298
$ws.Range("M98").Value = 498
$ws.Range("N98").Value = -4294

# row 113
$ws.Range("H113").Value = 71432310
$ws.Range("J113").Value = 8666.333000000001
$ws.Range("L113").Value = 8666.333000000001
$ws.Range("N113").Value = -15174.333

# row 116
$ws.Range("H116").Value = 4459.4
$ws.Range("I116").Value = 1860.625
$ws.Range("J116").Value = 7429.4287
$ws.Range("K116").Value = 1860.625
$ws.Range("L116").Value = 7429.4287
$ws.Range("M116").Value = 1581.375
$ws.Range("N116").Value = -14313.4287

# row 122
$ws.Range("H122").Value = 1059.6
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1298
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 3894
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -8794

$ws = $wb.Worksheets("ARM")
# row 2
$ws.Range("H2").Value = 2217.9443
$ws.Range("J2").Value = 4096.4
$ws.Range("L2").Value = 4096.4
$ws.Range("N2").Value = -4322.4

# row 4
$ws.Range("H4").Value = 59.75
$ws.Range("I4").Value = 64.666664
$ws.Range("K4").Value = 64.666664
$ws.Range("M4").Value = 51.333336

# row 61
$ws.Range("H61").Value = 3430.7693
$ws.Range("I61").Value = 1975
$ws.Range("J61").Value = 5760
$ws.Range("K61").Value = 1975
$ws.Range("L61").Value = 5760
$ws.Range("M61").Value = -1763
$ws.Range("N61").Value = -6184

# row 110
$ws.Range("H110").Value = 400
$ws.Range("I110").Value = 400
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 400
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1645
$ws.Range("N110").ClearContents()

# row 116
$ws.Range("H116").Value = 2217.9443
$ws.Range("J116").Value = 4096.4
$ws.Range("L116").Value = 4096.4
$ws.Range("N116").Value = -8684.4

# row 136
$ws.Range("H136").Value = 3430.7693
$ws.Range("I136").Value = 1975
$ws.Range("J136").Value = 5760
$ws.Range("K136").Value = 5925
$ws.Range("L136").Value = 17280
$ws.Range("M136").Value = -3375
$ws.Range("N136").Value = -22380

$ws = $wb.Worksheets("BSM")
# row 3
$ws.Range("H3").Value = 2217.9443
$ws.Range("J3").Value = 4096.4
$ws.Range("L3").Value = 4096.4
$ws.Range("N3").Value = -4324.4

# row 86
$ws.Range("H86").Value = 2100
$ws.Range("I86").Value = 1820
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 1820
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -697
$ws.Range("N86").Value = -5746

# row 89
$ws.Range("H89").Value = 2100
$ws.Range("I89").Value = 1820
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 9100
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -3484
$ws.Range("N89").Value = -28732

$ws = $wb.Worksheets("CUL")
# row 131
$ws.Range("H131").Value = 756.1900000000001
$ws.Range("J131").Value = 777.3579
$ws.Range("L131").Value = 2332.0737
$ws.Range("N131").Value = -12412.0737

$ws = $wb.Worksheets("GSM")
# row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# row 80
$ws.Range("H80").Value = 3801.2354
$ws.Range("I80").Value = 3480.125
$ws.Range("J80").Value = 4086.6667
$ws.Range("K80").Value = 3480.125
$ws.Range("L80").Value = 4086.6667
$ws.Range("M80").Value = -2482.125
$ws.Range("N80").Value = -6082.6667

# row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# row 83
$ws.Range("H83").Value = 3801.2354
$ws.Range("I83").Value = 3480.125
$ws.Range("J83").Value = 4086.6667
$ws.Range("K83").Value = 17400.625
$ws.Range("L83").Value = 20433.3335
$ws.Range("M83").Value = -12408.625
$ws.Range("N83").Value = -30417.3335

# row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# row 97
$ws.Range("H97").Value = 3500.111
$ws.Range("I97").Value = 1660.7142
$ws.Range("K97").Value = 1660.7142
$ws.Range("M97").Value = -1164.7142

# row 126
$ws.Range("H126").Value = 5117.241
$ws.Range("I126").Value = 3460
$ws.Range("J126").Value = 6892.857
$ws.Range("K126").Value = 10380
$ws.Range("L126").Value = 20678.571
$ws.Range("M126").Value = -7910
$ws.Range("N126").Value = -25618.571

$ws = $wb.Worksheets("LTW")
# row 7
$ws.Range("H7").Value = 3964.2856
$ws.Range("I7").Value = 4600
$ws.Range("J7").Value = 3116.6667
$ws.Range("K7").Value = 4600
$ws.Range("L7").Value = 3116.6667
$ws.Range("M7").Value = -4488
$ws.Range("N7").Value = -3340.6667

# row 40
$ws.Range("H40").Value = 5850
$ws.Range("I40").Value = 4575
$ws.Range("J40").Value = 7125
$ws.Range("K40").Value = 4575
$ws.Range("L40").Value = 7125
$ws.Range("M40").Value = -4439
$ws.Range("N40").Value = -7397

# row 46
$ws.Range("H46").Value = 2285
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 2670
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 2670
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -3046

# row 61
$ws.Range("H61").Value = 8061
$ws.Range("J61").Value = 10001.429
$ws.Range("L61").Value = 10001.429
$ws.Range("N61").Value = -10405.429

# row 93
$ws.Range("H93").Value = 2781.4614
$ws.Range("I93").Value = 2696.5833
$ws.Range("K93").Value = 2696.5833
$ws.Range("M93").Value = -1448.5833

# row 113
$ws.Range("H113").Value = 8061
$ws.Range("J113").Value = 10001.429
$ws.Range("L113").Value = 10001.429
$ws.Range("N113").Value = -14341.429

# row 122
$ws.Range("H122").Value = 1404295.6
$ws.Range("I122").Value = 3270749
$ws.Range("K122").Value = 9812247
$ws.Range("M122").Value = -9809797

# row 126
$ws.Range("H126").Value = 3964.2856
$ws.Range("I126").Value = 4600
$ws.Range("J126").Value = 3116.6667
$ws.Range("K126").Value = 13800
$ws.Range("L126").Value = 9350.000100000001
$ws.Range("M126").Value = -11330
$ws.Range("N126").Value = -14290.0001

$ws = $wb.Worksheets("WVR")
# row 62
$ws.Range("H62").Value = 2899.5
$ws.Range("I62").Value = 2899.5
$ws.Range("K62").Value = 2899.5
$ws.Range("M62").Value = -2275.5

# row 65
$ws.Range("H65").Value = 2899.5
$ws.Range("I65").Value = 2899.5
$ws.Range("K65").Value = 14497.5
$ws.Range("M65").Value = -11377.5

# row 107
$ws.Range("H107").Value = 4547067.5
$ws.Range("I107").Value = 825.8
$ws.Range("J107").Value = 9093309
$ws.Range("K107").Value = 2477.4
$ws.Range("L107").Value = 27279927
$ws.Range("M107").Value = -557.3999999999996
$ws.Range("N107").Value = -27283767

# row 126
$ws.Range("H126").Value = 2583.3333
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -13640
